$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new rows at row 66 (old rows 66-244 shift down to 68-246).
# ---------------------------------------------------------------------------
$ws.Rows("66:67").Insert()

# New row 66 stays blank, but needs the plain "blank row" style (s=6) rather
# than the default style Excel gives freshly inserted rows. Borrow the format
# from A3:D3 (a stable, untouched blank-style row above the insertion point).
$ws.Range("A3:D3").Copy()
$ws.Range("A66:D66").PasteSpecial(-4122)

# New row 67 becomes a category header ("주석"/"기타" pattern), matching the
# look of row 31 before it lost that label. Borrow formatting from row 31
# (still at its original position since it sits above the insertion point).
$ws.Range("A31:D31").Copy()
$ws.Range("A67:D67").PasteSpecial(-4122)
$ws.Range("A67").Value = "주석"
$ws.Range("C67").Value = "기타"

# ---------------------------------------------------------------------------
# 2. Row 31's category label changes from "기타" (Etc) to "Core".
#    Do this before introducing the other brand-new strings below so the
#    shared-string table gets "Core" appended first (index 187).
# ---------------------------------------------------------------------------
$ws.Range("C31").Value = "Core"

# ---------------------------------------------------------------------------
# 3. Row 71 (currently a blank filler row, style already s=6) receives the
#    new "incorrect_word" localization entry. Column order below matches the
#    original author's typing order so new shared strings land at indices
#    188, 189, 190 respectively.
# ---------------------------------------------------------------------------
$ws.Rows("71:71").RowHeight = 16
$ws.Range("D71").Value = "틀린단어"
$ws.Range("C71").Value = "incorrect word"
$ws.Range("A71").Value = "incorrect_word"

# ---------------------------------------------------------------------------
# 4. Add a thin border around the "incorrect word favorites" block (now at
#    rows 68-70) to match the refreshed styling for that section.
# ---------------------------------------------------------------------------
$ws.Range("A68:D70").Borders.LineStyle = 1

# ---------------------------------------------------------------------------
# 5. Two new blank rows (245, 246) are appended at the very end of the sheet,
#    matching the formatting of the last existing blank row (244).
# ---------------------------------------------------------------------------
$ws.Range("A244:D244").Copy()
$ws.Range("A245:D246").PasteSpecial(-4122)
$ws.Range("A245:D246").ClearContents()

# ---------------------------------------------------------------------------
# 6. Refresh the view: the sheet scrolled down and the selection moved from
#    A42 to A72.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 56
$ws.Range("A72").Select()
